$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells to be stored as text, so numeric-looking
# values (e.g. "1.00", "565.90") keep their exact original text
# representation instead of being coerced into floating point numbers.
# (Using one contiguous range -- a multi-area "A1,A2,..." range only applies
# the format to the first area in this runtime.)
$ws.Range("D2:D51").NumberFormat = "@"

# --- Row 14 / Row 15: WrappedliquidstakedEther2.0 and ShibaInu swapped places ---
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.891.94"
$ws.Range("E14").Value = "  -0.92%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "0.0000171"
$ws.Range("E15").Value = "  -4.68%  "

# --- Remaining Price / Volume(1h) updates ---
$ws.Range("D2").Value = "62.113.80"
$ws.Range("E2").Value = "  -1.90%  "
$ws.Range("D3").Value = "2.439.47"
$ws.Range("E3").Value = "  -1.54%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "565.90"
$ws.Range("E5").Value = "  -2.03%  "
$ws.Range("D6").Value = "144.82"
$ws.Range("E6").Value = "  -1.27%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "0.527"
$ws.Range("E8").Value = "  -2.81%  "
$ws.Range("D9").Value = "0.109"
$ws.Range("E9").Value = "  -2.61%  "
$ws.Range("E10").Value = "  -0.10%  "
$ws.Range("D11").Value = "5.19"
$ws.Range("E11").Value = "  -1.88%  "
$ws.Range("D12").Value = "0.345"
$ws.Range("E12").Value = "  -3.00%  "
$ws.Range("D13").Value = "28.46"
$ws.Range("E13").Value = "  -2.77%  "
$ws.Range("D16").Value = "62.146.89"
$ws.Range("E16").Value = "  -1.71%  "
$ws.Range("D17").Value = "2.446.09"
$ws.Range("E17").Value = "  -1.22%  "
$ws.Range("D18").Value = "7.68"
$ws.Range("E18").Value = "  -3.60%  "
$ws.Range("D19").Value = "10.67"
$ws.Range("E19").Value = "  -3.88%  "
$ws.Range("D20").Value = "318.95"
$ws.Range("E20").Value = "  -3.65%  "
$ws.Range("D21").Value = "4.09"
$ws.Range("E21").Value = "  -1.06%  "
$ws.Range("D22").Value = "2.14"
$ws.Range("E22").Value = "  -4.41%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").Value = "9.79"
$ws.Range("E24").Value = "  +6.32%  "
$ws.Range("D25").Value = "64.82"
$ws.Range("E25").Value = "  -2.40%  "
$ws.Range("D26").Value = "632.94"
$ws.Range("E26").Value = "  -5.61%  "
$ws.Range("D27").Value = "2.570.52"
$ws.Range("E27").Value = "  -1.04%  "
$ws.Range("D28").Value = "0.0₃0942"
$ws.Range("E28").Value = "  -6.73%  "
$ws.Range("D29").Value = "0.994"
$ws.Range("E29").Value = "  -0.42%  "
$ws.Range("D30").Value = "1.39"
$ws.Range("E30").Value = "  -5.31%  "
$ws.Range("D31").Value = "7.77"
$ws.Range("E31").Value = "  -5.02%  "
$ws.Range("D32").Value = "1.79"
$ws.Range("E32").Value = "  -4.58%  "
$ws.Range("D33").Value = "0.130"
$ws.Range("E33").Value = "  -5.61%  "
$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").Value = "1.46"
$ws.Range("E35").Value = "  -5.84%  "
$ws.Range("D36").Value = "4.60"
$ws.Range("E36").Value = "  -4.29%  "
$ws.Range("D37").Value = "150.19"
$ws.Range("E37").Value = "  -2.37%  "
$ws.Range("D38").Value = "0.362"
$ws.Range("E38").Value = "  -3.40%  "
$ws.Range("D39").Value = "18.33"
$ws.Range("E39").Value = "  -3.04%  "
$ws.Range("D40").Value = "5.19"
$ws.Range("E40").Value = "  -6.42%  "
$ws.Range("D41").Value = "2.67"
$ws.Range("E41").Value = "  -3.33%  "
$ws.Range("D42").Value = "1.68"
$ws.Range("E42").Value = "  -5.14%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").Value = "0.0₆0303"
$ws.Range("E44").Value = "  +0.98%  "
$ws.Range("D45").Value = "150.98"
$ws.Range("E45").Value = "  +1.69%  "
$ws.Range("D46").Value = "15.28"
$ws.Range("E46").Value = "  +0.84%  "
$ws.Range("D47").Value = "3.49"
$ws.Range("E47").Value = "  -3.82%  "
$ws.Range("D48").Value = "0.600"
$ws.Range("E48").Value = "  -1.38%  "
$ws.Range("D49").Value = "19.84"
$ws.Range("E49").Value = "  -5.21%  "
$ws.Range("D50").Value = "0.0498"
$ws.Range("E50").Value = "  -3.87%  "
$ws.Range("D51").Value = "0.0897"
$ws.Range("E51").Value = "  -2.85%  "
